$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Temps reel (j)" (actual time spent, column D) for existing tasks ---
$ws.Range("D3").Value = 14
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 2
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("D17").Value = 0

# --- Insert two new task rows (16 & 17) before the "Jalon" milestone rows,
#     pushing the milestone rows from 18-23 down to 20-25 ---
$ws.Rows("18:19").Insert()

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Création de fenêtres pop-up"
$ws.Range("D18").Value = 3

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Tests unitaires"
$ws.Range("D19").Value = 3

# --- Add the total row summing all actual-time entries ---
$ws.Range("D27").Formula = "=SUM(D3:D25)"

# --- Restore window/selection state ---
$ws.Range("K7").Select()
